$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 05:01"

# --- Row 37: Belgica (in-place value refresh) ---
$ws.Range("B37").Value = 66428
$ws.Range("C37").Value = 402
$ws.Range("D37").Value = 17452
$ws.Range("E37").Value = 39154
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 9822

# --- Rows 51-52: Honduras overtakes Barein ---
$ws.Range("A51").Value = "Honduras"
$ws.Range("B51").Value = 39741
$ws.Range("C51").Value = 465
$ws.Range("D51").Value = 5039
$ws.Range("E51").Value = 33536
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 50
$ws.Range("H51").Value = 1166

$ws.Range("A52").Value = "Barein"
$ws.Range("B52").Value = 39482
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 36110
$ws.Range("E52").Value = 3231
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 141

# --- Rows 70-75: Venezuela overtakes Costa Rica, Costa de Marfil, Chequia; Australia overtakes El Salvador ---
$ws.Range("A70").Value = "Venezuela"
$ws.Range("B70").Value = 15988
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 9959
$ws.Range("E70").Value = 5883
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 146

$ws.Range("A71").Value = "Costa Rica"
$ws.Range("B71").Value = 15841
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 3824
$ws.Range("E71").Value = 11902
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 115

$ws.Range("A72").Value = "Costa de Marfil"
$ws.Range("B72").Value = 15655
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 10361
$ws.Range("E72").Value = 5198
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 96

$ws.Range("A73").Value = "Chequia"
$ws.Range("B73").Value = 15516
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 11428
$ws.Range("E73").Value = 3715
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 373

$ws.Range("A74").Value = "Australia"
$ws.Range("B74").Value = 15302
$ws.Range("C74").Value = 367
$ws.Range("D74").Value = 9311
$ws.Range("E74").Value = 5824
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 6
$ws.Range("H74").Value = 167

$ws.Range("A75").Value = "El Salvador"
$ws.Range("B75").Value = 15035
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 7778
$ws.Range("E75").Value = 6849
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 408

# --- Row 91: Haiti (in-place value refresh) ---
$ws.Range("B91").Value = 7340
$ws.Range("C91").Value = 25
$ws.Range("D91").Value = 4365
$ws.Range("E91").Value = 2817
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 158
